$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.020.29"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +3.35%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.720.12"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("E4").Value = "  -0.03%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "218.67"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("E7").Value = "  +0.03%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "24.14"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +13.21%  "

$ws.Range("E9").Value = "  +2.95%  "

$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("E11").Value = "  +1.43%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.964.64"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.38%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.722.26"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.70%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.27"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.79%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.560"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +4.32%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "67.38"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.74%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "27.982.91"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.24%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "242.08"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.61%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0753"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.91"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.59%  "

$ws.Range("E21").Value = "  -0.07%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.61"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.32%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.66"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.77%  "

$ws.Range("E24").Value = "  +0.09%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "148.69"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("E26").Value = "  +3.62%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "16.67"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("E29").Value = "  -0.03%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0508"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.76%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.19"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.26"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.16%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.486.24"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -4.67%  "

$ws.Range("E35").Value = "  -2.91%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.952"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.47%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.606"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.42"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "

$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("E42").Value = "  +2.85%  "

$ws.Range("E43").Value = "  +0.00%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.29"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.868.54"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +2.22%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.796"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.75"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +11.17%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "90.81"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").Value = "  +3.52%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "8.25"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.16%  "

$ws.Range("E51").Value = "  +0.13%  "
